$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1310.3334
$ws.Range("I20").Value = 972.4
$ws.Range("J20").Value = 3000
$ws.Range("K20").Value = 972.4
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = -742.4
$ws.Range("N20").Value = -3460
$ws.Range("H35").Value = 1310.3334
$ws.Range("I35").Value = 972.4
$ws.Range("J35").Value = 3000
$ws.Range("K35").Value = 972.4
$ws.Range("L35").Value = 3000
$ws.Range("M35").Value = -593.4
$ws.Range("N35").Value = -3758
$ws.Range("H76").Value = 4080
$ws.Range("I76").Value = 3466.6667
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 3466.6667
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -3151.6667
$ws.Range("N76").Value = -5630
$ws.Range("H79").Value = 4080
$ws.Range("I79").Value = 3466.6667
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 3466.6667
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -2374.6667
$ws.Range("N79").Value = -7184
$ws.Range("H98").Value = 39624.234
$ws.Range("I98").Value = 55882.6
$ws.Range("J98").Value = 16398
$ws.Range("K98").Value = 55882.6
$ws.Range("L98").Value = 16398
$ws.Range("M98").Value = -54384.6
$ws.Range("N98").Value = -19394
$ws.Range("H112").Value = 1364.4
$ws.Range("J112").Value = 1476.4
$ws.Range("L112").Value = 4429.200000000001
$ws.Range("N112").Value = -6645.200000000001
$ws.Range("H122").Value = 39624.234
$ws.Range("I122").Value = 55882.6
$ws.Range("J122").Value = 16398
$ws.Range("K122").Value = 167647.8
$ws.Range("L122").Value = 49194
$ws.Range("M122").Value = -165197.8
$ws.Range("N122").Value = -54094
$ws.Range("H137").Value = 11826.728
$ws.Range("I137").Value = 25000.555
$ws.Range("K137").Value = 75001.66500000001
$ws.Range("M137").Value = -72451.66500000001
$ws.Range("H138").Value = 3852.8845
$ws.Range("J138").Value = 4786.325
$ws.Range("L138").Value = 14358.975
$ws.Range("N138").Value = -24638.975
$ws.Range("H141").Value = 4544.567
$ws.Range("I141").Value = 3761.6086
$ws.Range("J141").Value = 7117.143
$ws.Range("K141").Value = 11284.8258
$ws.Range("L141").Value = 21351.429
$ws.Range("M141").Value = -6104.825800000001
$ws.Range("N141").Value = -31711.429

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1482.5714
$ws.Range("I5").Value = 1315.8
$ws.Range("J5").Value = 1899.5
$ws.Range("K5").Value = 1315.8
$ws.Range("L5").Value = 1899.5
$ws.Range("M5").Value = -1203.8
$ws.Range("N5").Value = -2123.5
$ws.Range("H74").Value = 2623.5908
$ws.Range("I74").Value = 1412.6389
$ws.Range("K74").Value = 1412.6389
$ws.Range("M74").Value = -538.6388999999999
$ws.Range("H77").Value = 2623.5908
$ws.Range("I77").Value = 1412.6389
$ws.Range("K77").Value = 7063.1945
$ws.Range("M77").Value = -2695.1945

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1482.5714
$ws.Range("I4").Value = 1315.8
$ws.Range("J4").Value = 1899.5
$ws.Range("K4").Value = 1315.8
$ws.Range("L4").Value = 1899.5
$ws.Range("M4").Value = -1200.8
$ws.Range("N4").Value = -2129.5
$ws.Range("H99").Value = 15212.857
$ws.Range("I99").Value = 16821.25
$ws.Range("J99").Value = 5562.5
$ws.Range("K99").Value = 16821.25
$ws.Range("L99").Value = 5562.5
$ws.Range("M99").Value = -15323.25
$ws.Range("N99").Value = -8558.5
$ws.Range("H105").Value = 3559.6155
$ws.Range("I105").Value = 2697.2222
$ws.Range("J105").Value = 5500
$ws.Range("K105").Value = 2697.2222
$ws.Range("L105").Value = 5500
$ws.Range("M105").Value = -950.2222000000002
$ws.Range("N105").Value = -8994
$ws.Range("H134").Value = 2752.276
$ws.Range("I134").Value = 1906.6086
$ws.Range("K134").Value = 5719.825800000001
$ws.Range("M134").Value = -3184.825800000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5598.3335
$ws.Range("I31").Value = 1295
$ws.Range("K31").Value = 1295
$ws.Range("M31").Value = -1000
$ws.Range("H34").Value = 5598.3335
$ws.Range("I34").Value = 1295
$ws.Range("K34").Value = 1295
$ws.Range("M34").Value = -1093
$ws.Range("H86").Value = 12260.107
$ws.Range("I86").Value = 11573
$ws.Range("K86").Value = 11573
$ws.Range("M86").Value = -10450
$ws.Range("H89").Value = 12260.107
$ws.Range("I89").Value = 11573
$ws.Range("K89").Value = 57865
$ws.Range("M89").Value = -52249

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 81
$ws.Range("I7").Value = 81
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 243
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -131
$ws.Range("H38").Value = 1240.963
$ws.Range("I38").Value = 305.69232
$ws.Range("J38").Value = 2109.4285
$ws.Range("K38").Value = 917.07696
$ws.Range("L38").Value = 6328.2855
$ws.Range("M38").Value = -570.07696
$ws.Range("N38").Value = -7022.2855
$ws.Range("H51").Value = 1757.1428
$ws.Range("I51").Value = 975.7
$ws.Range("K51").Value = 2927.1
$ws.Range("M51").Value = -2467.1
$ws.Range("H139").Value = 1501759.4
$ws.Range("I139").Value = 1668066
$ws.Range("K139").Value = 5004198
$ws.Range("M139").Value = -4999058

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 84.40909000000001
$ws.Range("I2").Value = 84.40909000000001
$ws.Range("K2").Value = 84.40909000000001
$ws.Range("M2").Value = 28.59090999999999
$ws.Range("H102").Value = 18887.375
$ws.Range("I102").Value = 24049.834
$ws.Range("J102").Value = 3400
$ws.Range("K102").Value = 24049.834
$ws.Range("L102").Value = 3400
$ws.Range("M102").Value = -22427.834
$ws.Range("N102").Value = -6644
$ws.Range("H123").Value = 21899.8
$ws.Range("J123").Value = 21899.8
$ws.Range("L123").Value = 21899.8
$ws.Range("N123").Value = -26799.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 38400.6
$ws.Range("J38").Value = 42493.25
$ws.Range("L38").Value = 42493.25
$ws.Range("N38").Value = -43313.25
$ws.Range("H47").Value = 49000
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H52").Value = 49000
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H122").Value = 4313.5186
$ws.Range("I122").Value = 3941.875
$ws.Range("J122").Value = 4854.091
$ws.Range("K122").Value = 11825.625
$ws.Range("L122").Value = 14562.273
$ws.Range("M122").Value = -9375.625
$ws.Range("N122").Value = -19462.273
$ws.Range("H141").Value = 583107.5
$ws.Range("J141").Value = 583107.5
$ws.Range("L141").Value = 583107.5
$ws.Range("N141").Value = -593467.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 99999
$ws.Range("J46").Value = 99999
$ws.Range("L46").Value = 99999
$ws.Range("N46").Value = -100461
$ws.Range("H64").Value = 94000
$ws.Range("J64").Value = 94000
$ws.Range("L64").Value = 94000
$ws.Range("N64").Value = -94496
$ws.Range("H67").Value = 94000
$ws.Range("J67").Value = 94000
$ws.Range("L67").Value = 94000
$ws.Range("N67").Value = -95716
$ws.Range("H81").Value = 15384.818
$ws.Range("I81").Value = 21700
$ws.Range("J81").Value = 4333.25
$ws.Range("K81").Value = 43400
$ws.Range("L81").Value = 8666.5
$ws.Range("M81").Value = -42339
$ws.Range("N81").Value = -10788.5
$ws.Range("H84").Value = 15384.818
$ws.Range("I84").Value = 21700
$ws.Range("J84").Value = 4333.25
$ws.Range("K84").Value = 217000
$ws.Range("L84").Value = 43332.5
$ws.Range("M84").Value = -211696
$ws.Range("N84").Value = -53940.5
$ws.Range("H96").Value = 2070.75
$ws.Range("I96").Value = 1468.7
$ws.Range("K96").Value = 1468.7
$ws.Range("M96").Value = -95.70000000000005
$ws.Range("H134").Value = 99999
$ws.Range("J134").Value = 99999
$ws.Range("L134").Value = 299997
$ws.Range("N134").Value = -305067
